# ---------------------------------------------------------------------------
# Edit script for "7) Contact Page.docx"
#
#  1) Paragraph "anne_zambrano@dlsu.edu.ph" (2nd-level bullet): change its
#     spacing from w:after="0" w:afterAutospacing="0" to w:after="240"
#     (drop the afterAutospacing attribute).
#  2) Collapse the three paragraphs "Social links" / (empty 2nd-level bullet)
#     / "Disclaimer" into a single, non-bulleted paragraph: drop the numbering,
#     bump the run size to 30/30, change the indent to left=720/firstLine=0,
#     keep spacing after=240, and remove all the paragraphs' text runs except
#     for the final empty trailing run.
#  3) Style tweaks: add the "complex script" (Cs) companion toggle next to the
#     existing b/i toggles in Heading 2, Heading 3, Heading 6 and Subtitle.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Get-ParagraphIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# --- 1) Fix spacing on the "anne_zambrano@dlsu.edu.ph" paragraph ---------
$idxAnne = Get-ParagraphIndexByText $d "anne_zambrano@dlsu.edu.ph"
$pAnne = $d.Paragraphs.Item($idxAnne)
$xmlAnne = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="240" w:lineRule="auto"/><w:ind w:left="1440" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Roboto Serif" w:cs="Roboto Serif" w:eastAsia="Roboto Serif" w:hAnsi="Roboto Serif"/></w:rPr></w:pPr><w:hyperlink r:id="rId7"><w:r><w:rPr><w:rFonts w:ascii="Roboto Serif" w:cs="Roboto Serif" w:eastAsia="Roboto Serif" w:hAnsi="Roboto Serif"/><w:color w:val="1155cc"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">anne_zambrano@dlsu.edu.ph</w:t></w:r></w:hyperlink><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $pAnne.Range.InsertXML($xmlAnne)

# --- 2) Merge "Social links" / empty / "Disclaimer" paragraphs into one --
$idxSocial = Get-ParagraphIndexByText $d "Social links"
$idxDisclaimer = Get-ParagraphIndexByText $d "Disclaimer"
$pSocial = $d.Paragraphs.Item($idxSocial)
$pDisclaimer = $d.Paragraphs.Item($idxDisclaimer)
$mergeRange = $d.Range($pSocial.Range.Start, $pDisclaimer.Range.End)
$xmlMerge = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="240" w:lineRule="auto"/><w:ind w:left="720" w:firstLine="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Roboto Serif" w:cs="Roboto Serif" w:eastAsia="Roboto Serif" w:hAnsi="Roboto Serif"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$mergeRange.InsertXML($xmlMerge)

# InsertXML canonicalizes away indent values that equal Word's implicit
# default (firstLine=0), so set it explicitly through the object model too.
$pMerged = $d.Paragraphs.Item($idxSocial)
$pMerged.Format.LeftIndent = 36
$pMerged.Format.FirstLineIndent = 0

# --- 3) Style updates: add complex-script (Cs) counterparts --------------
$heading2 = $d.Styles.Item("Heading 2")
$heading2.Font.BoldBi = 0

$heading3 = $d.Styles.Item("Heading 3")
$heading3.Font.BoldBi = 0

$heading6 = $d.Styles.Item("Heading 6")
$heading6.Font.ItalicBi = 1

$subtitle = $d.Styles.Item("Subtitle")
$subtitle.Font.ItalicBi = 0

Write-Host "Edits applied"
